# OrangeHRMUserRoles.xlsx - add the ESS_users sheet (list of ESS usernames)
# after Sheet1, and restore the A2 selection on Sheet1.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)

# New worksheet, inserted right after Sheet1.
$essUsers = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$essUsers.Name = "ESS_users"

$users = @(
    "Aaliyah.Haq",
    "Aatmaram",
    "Alice.Duval",
    "Anthony.Nolan",
    "Aravind",
    "azq@gmail.com",
    "Cassidy.Hope",
    "catelusCuParulCret",
    "Cecil.Bonaparte",
    "Charlie.Carter",
    "Chenzira.Chuki",
    "David.Morris",
    "Ehioze.Ebo",
    "Fiona.Grace",
    "Garry.White",
    "Goutam.Ganesh",
    "Jacqueline.White",
    "Jadine.Jackie",
    "Jasmine.Morgan",
    "Joe.Root",
    "Jordan.Mathews",
    "Kevin.Mathews",
    "Kiyara.Hu",
    "Lisa.Andrews",
    "Luke.Wright",
    "Maggie.Manning",
    "manali28",
    "Melan.Peiris",
    "Nathan.Elliot",
    "Nina.Patel",
    "Rebecca.Harmony",
    "Russel.Hamilton",
    "Sania.Shaheen",
    "Sara.Tencrady"
)

for ($i = 0; $i -lt $users.Length; $i++) {
    $row = $i + 1
    $essUsers.Cells.Item($row, 1).Value = $users[$i]
}

# Restore Sheet1 as the active sheet with A2 selected.
$sheet1.Activate()
$sheet1.Range("A2").Select() | Out-Null
